$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New stress results rows to append (rows 57-59)
$newRows = @(
    @("2023-07-13 21:57:32", "Cifar10", 100, 10, 2, 8, 0.1302874556367354),
    @("2023-07-13 21:57:41", "Cifar10", 100, 30, 8, 22, 0.05650222861550116),
    @("2023-07-13 21:59:47", "Cifar10", 100, 100, 34, 66, 0.0181339383423586)
)

$startRow = 57
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
